# Insert a new data row at row 51 (pushes existing rows 51..111 down to 52..112)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51:51").Insert()

$newRow = 51

$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44792
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112052
$ws.Cells.Item($newRow, 7).Value = "Albahaca"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 1000
$ws.Cells.Item($newRow, 11).Value = 3400
$ws.Cells.Item($newRow, 12).Value = 3500
$ws.Cells.Item($newRow, 13).Value = 3450
$ws.Cells.Item($newRow, 14).Value = "`$/paquete"
$ws.Cells.Item($newRow, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($newRow, 16).Value = 3450
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
